# Generate Report for Handoff
#
# The localization-status report has moved from "In Translation" to
# "Ready for handoff": update the Status columns and the corresponding
# "Latest Handoff"/"Latest HO Xliff Generate" timestamps on each sheet,
# then resize the Status columns so the longer text is fully visible.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status: "In Translation" -> "Ready for handoff" ------------------
$wsOverview.Range("E2").Value = "Ready for handoff"   # zh-cn status
$wsOverview.Range("F2").Value = "Ready for handoff"   # de-de status
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("C2").Value = "Ready for handoff"

# --- Handoff timestamps, bumped forward by 30 seconds ------------------
$wsOverview.Range("G2").Value = "2016-08-15 12:53:46" # Latest HO Xliff Generate Date
$wsDeDe.Range("H2").Value     = "2016-08-15 12:53:46" # de-de Latest Handoff Datetime
$wsZhCn.Range("H2").Value     = "2016-08-15 12:53:41" # zh-cn Latest Handoff Datetime

# --- Widen the Status columns to fit "Ready for handoff" ---------------
$statusColumnWidth = 16.333333333333332
$wsOverview.Range("E1:F1").EntireColumn.ColumnWidth = $statusColumnWidth
$wsZhCn.Range("C1").EntireColumn.ColumnWidth = $statusColumnWidth
$wsDeDe.Range("C1").EntireColumn.ColumnWidth = $statusColumnWidth
